# Auto-generated edit script applying the "CI: Auto Update Data" diff
# to 悖论模拟干员名单作者版.xlsx (single worksheet, Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.2), maa://25390 (95.85), maa://36681 (86.67)"
$ws.Range("L2").Value = "*maa://24633 (55.26), *maa://30515 (69.0), *maa://34787 (71.88), ***maa://20792 (11.93), maa://39402 (84.85), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.72), *maa://20791 (62.32)"
$ws.Range("AB2").Value = "maa://21246 (91.32), maa://36684 (97.65), ***maa://22731 (6.67)"
$ws.Range("AF2").Value = "maa://25251 (92.31), ***maa://21730 (18.18), ***maa://39501 (15.79), *maa://36675 (60.0)"
$ws.Range("X3").Value = "maa://27396 (84.54), maa://27484 (95.92), maa://27480 (82.35)"
$ws.Range("D4").Value = "maa://24632 (93.79), **maa://24303 (33.33), maa://22499 (85.71), maa://22746 (100.0)"
$ws.Range("X4").Value = "**maa://32495 (47.27), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (88.24)"
$ws.Range("D5").Value = "maa://21245 (82.69), maa://22744 (84.0)"
$ws.Range("H6").Value = "maa://24370 (96.49)"
$ws.Range("P6").Value = "maa://31836 (90.0), maa://30381 (92.31)"
$ws.Range("A8").Value = "更新日期：2024.11.24 12:07:34"
$ws.Range("X9").Value = "maa://26223 (97.39)"
$ws.Range("T10").Value = "maa://27395 (95.88), maa://22755 (87.39), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("T11").Value = "maa://22747 (93.2), maa://22501 (98.36)"
$ws.Range("X11").Value = "maa://36713 (98.1)"
$ws.Range("D13").Value = "maa://24999 (91.63), maa://36673 (92.42), maa://25001 (85.51)"
$ws.Range("L14").Value = "maa://26245 (96.3), maa://21288 (96.21), maa://36682 (97.37), maa://39841 (94.12)"
$ws.Range("AB14").Value = "maa://22764 (96.77)"
$ws.Range("D16").Value = "maa://21441 (96.28), maa://36679 (92.68), maa://37650 (96.77)"
$ws.Range("P16").Value = "maa://28504 (90.74)"
$ws.Range("T16").Value = "maa://22729 (95.3), *maa://28648 (67.24), maa://36674 (82.05)"
$ws.Range("X16").Value = "maa://28501 (97.8), maa://28051 (96.0)"
$ws.Range("T18").Value = "maa://24385 (97.06)"
$ws.Range("AB18").Value = "maa://24393 (97.37)"
$ws.Range("T19").Value = "maa://24386 (98.95)"
$ws.Range("D20").Value = "maa://21432 (90.21), maa://25198 (92.93), *maa://20795 (50.79), maa://36680 (96.43)"
$ws.Range("L20").Value = "maa://41331 (81.93)"
$ws.Range("P21").Value = "maa://24381 (87.5)"
$ws.Range("T21").Value = "maa://21993 (88.89)"
$ws.Range("L23").Value = "maa://39756 (93.21), maa://39875 (93.22)"
$ws.Range("X24").Value = "maa://29988 (86.36), maa://23504 (93.07), **maa://22892 (39.86), *maa://25141 (77.6), maa://36663 (80.95), ***maa://22815 (23.08)"
$ws.Range("D25").Value = "maa://29753 (95.12)"
$ws.Range("P25").Value = "maa://24382 (93.1)"
$ws.Range("AB25").Value = "maa://31215 (85.11), *maa://24516 (79.07), maa://26001 (87.27)"
$ws.Range("H26").Value = "maa://24913 (91.25)"
$ws.Range("H27").Value = "**maa://21283 (48.65), maa://34494 (96.3), *maa://39601 (73.33), **maa://36665 (44.44)"
$ws.Range("T27").Value = "*maa://30624 (75.0)"
$ws.Range("X28").Value = "maa://39929 (89.35), ***maa://39723 (14.29), maa://41749 (86.84)"
$ws.Range("AF28").Value = "maa://36660 (92.59), *maa://36701 (62.96)"
$ws.Range("H29").Value = "*maa://25175 (69.57)"
$ws.Range("P29").Value = "*maa://23168 (54.72), **maa://30050 (43.48)"
$ws.Range("AF29").Value = "*maa://24080 (69.23), ***maa://34960 (8.33), maa://42865 (82.14)"
$ws.Range("AB30").Value = "maa://42979 (97.26)"
$ws.Range("L31").Value = "maa://35926 (93.82), maa://36258 (81.82)"
$ws.Range("T32").Value = "maa://41108 (87.76), maa://42859 (93.88), maa://41238 (95.38)"
$ws.Range("L35").Value = "maa://41296 (95.92)"
$ws.Range("H37").Value = "*maa://24374 (58.14)"
$ws.Range("T37").Value = "**maa://39354 (33.33)"
$ws.Range("AF38").Value = "maa://36697 (86.06)"
$ws.Range("H39").Value = "maa://25199 (85.32), maa://36670 (88.46), maa://30434 (88.52), ***maa://25036 (16.0)"
$ws.Range("P39").Value = "maa://24709 (91.53)"
# G46: value is purely numeric-looking ("2") but must stay text,
# like the original cell. Route it through a scratch cell so Excel
# does not auto-convert it to a number when writing.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "2"
$scratch.Copy()
$ws.Range("G46").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("H46").Value = "maa://35931 (92.57), maa://43901 (100.0)"
$ws.Range("H47").Value = "maa://27410 (96.05), maa://29661 (97.78), maa://28038 (84.62)"
$ws.Range("H52").Value = "maa://24376 (96.55)"
$ws.Range("H55").Value = "maa://32532 (92.37)"
$ws.Range("H57").Value = "maa://25176 (98.15)"
$ws.Range("H60").Value = "*maa://40438 (56.41)"
$ws.Range("H62").Value = "maa://42981 (95.45)"
